# CoastSnapPyDB.xlsx: introduce an ACTIVE/INACTIVE marker in cell B1 of
# each station sheet and make "texel" the active tab/sheet (was "egmond").

$wb = $excel.ActiveWorkbook

$egmond = $wb.Worksheets.Item("egmond")
$texel  = $wb.Worksheets.Item("texel")

# egmond is the currently-active station -> ACTIVE
$egmond.Range("B1").Value = "ACTIVE"

# texel is not currently active -> INACTIVE
$texel.Range("B1").Value = "INACTIVE"

# Reset the lingering selections on both sheets back to B1 (top-left),
# then make "texel" the selected / active tab, as in the saved file.
$egmond.Activate()
$egmond.Range("B1").Select()

$texel.Activate()
$texel.Range("B1").Select()
